$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (apptyp_code) holds numeric codes like "001".."016" that now
# need an "A" prefix (e.g. "A001".."A016"). Walk the used rows and patch
# any matching values.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 64 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null) {
        $s = [string]$val
        if ($s -match '^[0-9]{3}$') {
            $cell.Value2 = "A" + $s
        }
    }
}

# Reflect the new selection / scroll position recorded in the workbook view.
$ws.Range("M55").Select()
$excel.ActiveWindow.ScrollRow = 43
